$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 400-401 (pushes existing rows 400-422 down to 402-424),
# inheriting formatting from the row above (row 399), matching the diff's
# dimension change from A1:T422 to A1:T424.
$ws.Range("A400:T401").EntireRow.Insert()

# New row 400: Frutilla, Especial, most recent week's data point
$ws.Range("A400").Value = 7
$ws.Range("B400").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C400").Value = "Ñuble"
$ws.Range("D400").Value = 44931
$ws.Range("E400").Value = 16
$ws.Range("F400").Value = "Fruta"
$ws.Range("G400").Value = 100101
$ws.Range("H400").Value = "Berries"
$ws.Range("I400").Value = 100112025
$ws.Range("J400").Value = "Frutilla"
$ws.Range("K400").Value = "Sin especificar"
$ws.Range("L400").Value = "Especial"
$ws.Range("M400").Value = 80
$ws.Range("N400").Value = 8000
$ws.Range("O400").Value = 8000
$ws.Range("P400").Value = 8000
$ws.Range("Q400").Value = "$/caja 7 kilos"
$ws.Range("R400").Value = "Provincia de Diguillín"
$ws.Range("S400").Value = 1143
$ws.Range("T400").Value = 7

# New row 401: Frutilla, Primera, most recent week's data point
$ws.Range("A401").Value = 7
$ws.Range("B401").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C401").Value = "Ñuble"
$ws.Range("D401").Value = 44931
$ws.Range("E401").Value = 16
$ws.Range("F401").Value = "Fruta"
$ws.Range("G401").Value = 100101
$ws.Range("H401").Value = "Berries"
$ws.Range("I401").Value = 100112025
$ws.Range("J401").Value = "Frutilla"
$ws.Range("K401").Value = "Sin especificar"
$ws.Range("L401").Value = "Primera"
$ws.Range("M401").Value = 120
$ws.Range("N401").Value = 6500
$ws.Range("O401").Value = 7000
$ws.Range("P401").Value = 6750
$ws.Range("Q401").Value = "$/caja 7 kilos"
$ws.Range("R401").Value = "Provincia de Diguillín"
$ws.Range("S401").Value = 964
$ws.Range("T401").Value = 7
